$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.847.93"
$ws.Range("D3").Value = "2.624.75"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "520.63"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "145.14"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "2.638.94"
$ws.Range("E9").Value = "  +2.45%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "6.30"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "3.088.44"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "58.865.56"
$ws.Range("E15").Value = "  +0.75%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "20.85"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "2.636.06"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E23").Value = "  -0.16%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "61.54"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  -0.25%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.165"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "0.0₃0799"
$ws.Range("E28").Value = "  -1.32%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.10"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +3.17%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "18.85"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  +2.01%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "150.46"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.48%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.978"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  +1.09%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "36.63"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  +1.53%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "277.60"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("E45").Value = "  -0.25%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "19.48"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("E47").Value = "  -2.41%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "10.30"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "1.988.73"
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -0.79%  "
